$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 5899
$ws1.Range("F7").Value = 540
$ws1.Range("F9").Value = 1571
$ws1.Range("F10").Value = 19
$ws1.Range("F11").Value = 31
$ws1.Range("F12").Value = 676
$ws1.Range("F15").Value = 1564
$ws1.Range("F16").Value = 553
$ws1.Range("F17").Value = 152
$ws1.Range("F18").Value = 625
$ws1.Range("F19").Value = 4447
$ws1.Range("F22").Value = 3340
$ws1.Range("F23").Value = 816
$ws1.Range("F24").Value = 13
$ws1.Range("F25").Value = 48
$ws1.Range("F26").Value = 2309
$ws1.Range("F28").Value = 336
$ws1.Range("F30").Value = 453
$ws1.Range("F31").Value = 1226
$ws1.Range("F32").Value = 787
$ws1.Range("F33").Value = 2
$ws1.Range("F34").Value = 3
$ws1.Range("F36").Value = 1211
$ws1.Range("F37").Value = 1195

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 106
$ws2.Range("F15").Value = 48

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 680
$ws3.Range("F4").Value = 185
$ws3.Range("F5").Value = 272

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 680
$ws4.Range("F7").Value = 185
$ws4.Range("F8").Value = 5899
$ws4.Range("F13").Value = 106
$ws4.Range("F16").Value = 540
$ws4.Range("F19").Value = 1571
$ws4.Range("F21").Value = 19
$ws4.Range("F22").Value = 31
$ws4.Range("F25").Value = 1564
$ws4.Range("F26").Value = 553
$ws4.Range("F27").Value = 152
$ws4.Range("F28").Value = 625
$ws4.Range("F29").Value = 4447
$ws4.Range("F31").Value = 3340
$ws4.Range("F32").Value = 816
$ws4.Range("F33").Value = 48
$ws4.Range("F35").Value = 2309
$ws4.Range("F37").Value = 336
$ws4.Range("F39").Value = 453
$ws4.Range("F40").Value = 1226
$ws4.Range("F45").Value = 787
$ws4.Range("F47").Value = 1211
$ws4.Range("F49").Value = 1195
